$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the floating point precision of the existing A5 timestamp
$ws.Range("A5").Value = 45806.39326444444

# Append the new row of price data
$ws.Range("A6").Value = 45806.4066544568
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat

$ws.Range("B6").Value = "EVOWHEY PROTEIN"
$ws.Range("C6").Value = "2Kg"
$ws.Range("D6").Value = "37,90€"
